$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.09260899626369
$ws.Range("C2").Value = 11.01616173040622
$ws.Range("E2").Value = 8.830109479115711
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.735465500740413
$ws.Range("I2").Value = 37.6181221768755
$ws.Range("L2").Value = 10.49362093115385
$ws.Range("M2").Value = 17.61542889731175
$ws.Range("B3").Value = 19.83923396387046
$ws.Range("C3").Value = 10.49671138968456
$ws.Range("E3").Value = 8.799417286426046
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.739965927723149
$ws.Range("I3").Value = 37.27220863289495
$ws.Range("L3").Value = 10.50252392224762
$ws.Range("M3").Value = 17.58909749825791
$ws.Range("B4").Value = 19.68955781947197
$ws.Range("C4").Value = 10.16803367322684
$ws.Range("E4").Value = 8.780070512100735
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.742867962003373
$ws.Range("I4").Value = 37.06322622980638
$ws.Range("L4").Value = 10.50956765205805
$ws.Range("M4").Value = 17.57734507447164
$ws.Range("B5").Value = 19.63012639586855
$ws.Range("C5").Value = 10.03185764914283
$ws.Range("E5").Value = 8.77206012604975
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.744085604945199
$ws.Range("I5").Value = 36.97896853196848
$ws.Range("L5").Value = 10.51283440494702
$ws.Range("M5").Value = 17.57366902273553
$ws.Range("B6").Value = 19.62035450758301
$ws.Range("C6").Value = 10.00911732876797
$ws.Range("E6").Value = 8.770722364079907
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.744289914367977
$ws.Range("I6").Value = 36.96503354642433
$ws.Range("L6").Value = 10.51340077842698
$ws.Range("M6").Value = 17.57312591656699
$ws.Range("B7").Value = 19.68874987574315
$ws.Range("C7").Value = 10.166205914505
$ws.Range("E7").Value = 8.779962993470464
$ws.Range("F7").Value = 15.26647399323728
$ws.Range("G7").Value = 3.742884241484834
$ws.Range("I7").Value = 37.06208617929367
$ws.Range("L7").Value = 10.50961010407838
$ws.Range("M7").Value = 17.57729098760325
$ws.Range("B8").Value = 20.00407035906811
$ws.Range("C8").Value = 10.83919259636087
$ws.Range("E8").Value = 8.819630756316489
$ws.Range("F8").Value = 16.53996406344769
$ws.Range("G8").Value = 3.736988539912854
$ws.Range("I8").Value = 37.4981722764414
$ws.Range("L8").Value = 10.49636319924594
$ws.Range("M8").Value = 17.60543548958041
$ws.Range("B9").Value = 20.66553261144694
$ws.Range("C9").Value = 12.07387022820516
$ws.Range("E9").Value = 8.893453436283853
$ws.Range("F9").Value = 19.0027458068253
$ws.Range("G9").Value = 3.726521200743181
$ws.Range("I9").Value = 38.37823640381075
$ws.Range("L9").Value = 10.48290983844246
$ws.Range("M9").Value = 17.69549182920973
$ws.Range("B10").Value = 21.1727189650717
$ws.Range("C10").Value = 12.92044394064993
$ws.Range("E10").Value = 8.945296887830215
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.719488260937742
$ws.Range("I10").Value = 39.03699851607502
$ws.Range("L10").Value = 10.48067020368067
$ws.Range("M10").Value = 17.78263097808993
$ws.Range("B11").Value = 21.40695567827485
$ws.Range("C11").Value = 13.29103059163995
$ws.Range("E11").Value = 8.968367180049718
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.716429488221879
$ws.Range("I11").Value = 39.33866255763698
$ws.Range("L11").Value = 10.48131177952981
$ws.Range("M11").Value = 17.82674875433409
$ws.Range("B12").Value = 21.49606866397237
$ws.Range("C12").Value = 13.42918181533485
$ws.Range("E12").Value = 8.97702978409643
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.715291261578985
$ws.Range("I12").Value = 39.45312291267609
$ws.Range("L12").Value = 10.48179334059945
$ws.Range("M12").Value = 17.84409077027777
$ws.Range("B13").Value = 21.47685970152366
$ws.Range("C13").Value = 13.39952683878168
$ws.Range("E13").Value = 8.975167407634901
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.715535509018491
$ws.Range("I13").Value = 39.42846269454414
$ws.Range("L13").Value = 10.48167901895526
$ws.Range("M13").Value = 17.84032773133115
$ws.Range("B14").Value = 21.4142793150904
$ws.Range("C14").Value = 13.3024405835158
$ws.Range("E14").Value = 8.96908133170114
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.716335444352195
$ws.Range("I14").Value = 39.34807507802368
$ws.Range("L14").Value = 10.48134661694406
$ws.Range("M14").Value = 17.82816280537411
$ws.Range("B15").Value = 21.37599802224105
$ws.Range("C15").Value = 13.24268579929005
$ws.Range("E15").Value = 8.965343854602891
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.716828036452432
$ws.Range("I15").Value = 39.29886309657738
$ws.Range("L15").Value = 10.48117407911757
$ws.Range("M15").Value = 17.82079394533515
$ws.Range("B16").Value = 21.15747381673244
$ws.Range("C16").Value = 12.89592462467951
$ws.Range("E16").Value = 8.94377887587363
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.719690974585222
$ws.Range("I16").Value = 39.01731903694032
$ws.Range("L16").Value = 10.48066167563323
$ws.Range("M16").Value = 17.77983718978147
$ws.Range("B17").Value = 21.02425022148241
$ws.Range("C17").Value = 12.67940634493011
$ws.Range("E17").Value = 8.930418091979146
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.721483187721321
$ws.Range("I17").Value = 38.84506844783681
$ws.Range("L17").Value = 10.48077252025998
$ws.Range("M17").Value = 17.75585276078988
$ws.Range("B18").Value = 20.94795990265688
$ws.Range("C18").Value = 12.55350817412092
$ws.Range("E18").Value = 8.922684970262138
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.722527258994921
$ws.Range("I18").Value = 38.74618536196975
$ws.Range("L18").Value = 10.48099257451826
$ws.Range("M18").Value = 17.74247957243249
$ws.Range("B19").Value = 20.9221899291199
$ws.Range("C19").Value = 12.51065030111818
$ws.Range("E19").Value = 8.920058356530458
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.722883041628049
$ws.Range("I19").Value = 38.71273983535781
$ws.Range("L19").Value = 10.48109392882726
$ws.Range("M19").Value = 17.73802436228596
$ws.Range("B20").Value = 21.03839796373638
$ws.Range("C20").Value = 12.70259683444312
$ws.Range("E20").Value = 8.93184537281436
$ws.Range("F20").Value = 20.2495528364879
$ws.Range("G20").Value = 3.721291034583776
$ws.Range("I20").Value = 38.86338548857049
$ws.Range("L20").Value = 10.48074454509242
$ws.Range("M20").Value = 17.75836232417328
$ws.Range("B21").Value = 21.43265022689819
$ws.Range("C21").Value = 13.33101705256678
$ws.Range("E21").Value = 8.970870955749662
$ws.Range("F21").Value = 21.46857628470567
$ws.Range("G21").Value = 3.716099940582411
$ws.Range("I21").Value = 39.37168119219635
$ws.Range("L21").Value = 10.48143777732869
$ws.Range("M21").Value = 17.831718760393
$ws.Range("B22").Value = 21.69267990249081
$ws.Range("C22").Value = 13.72897554938069
$ws.Range("E22").Value = 8.995947456493916
$ws.Range("F22").Value = 22.22866616901555
$ws.Range("G22").Value = 3.712824142238319
$ws.Range("I22").Value = 39.70517552857493
$ws.Range("L22").Value = 10.48328151233076
$ws.Range("M22").Value = 17.88336156526955
$ws.Range("B23").Value = 21.55371128107053
$ws.Range("C23").Value = 13.51777116291065
$ws.Range("E23").Value = 8.98260282941337
$ws.Range("F23").Value = 21.82633154475864
$ws.Range("G23").Value = 3.714561851093781
$ws.Range("I23").Value = 39.52708469132782
$ws.Range("L23").Value = 10.48217030409874
$ws.Range("M23").Value = 17.85546323729773
$ws.Range("B24").Value = 21.03200082267778
$ws.Range("C24").Value = 12.69211682653014
$ws.Range("E24").Value = 8.931200260602257
$ws.Range("F24").Value = 20.22900810905294
$ws.Range("G24").Value = 3.721377864321424
$ws.Range("I24").Value = 38.85510389497359
$ws.Range("L24").Value = 10.48075670572753
$ws.Range("M24").Value = 17.75722645488716
$ws.Range("B25").Value = 20.48249178538151
$ws.Range("C25").Value = 11.74994990126773
$ws.Range("E25").Value = 8.873904561889107
$ws.Range("F25").Value = 18.34778573295697
$ws.Range("G25").Value = 3.729236752506505
$ws.Range("I25").Value = 38.137793264595
$ws.Range("L25").Value = 10.48520723212248
$ws.Range("M25").Value = 17.66742187987323
